# Weekly refresh: a new "Sandia" price record for Agrícola del Norte S.A. de
# Arica is inserted ahead of the existing log (new row 40), pushing the
# previously-logged rows 40-61 down to 41-62. The sheet's other columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Unidad de comercialización, Kg o Unidades, Clasificación) are identical for
# every row in this block, so only the row-shift + the new row's values are
# needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 40, shifting existing rows 40:61 down to 41:62
# (xlShiftDown = -4121). Excel copies formatting from the row above, which
# already carries the date style (s="2") on column D.
$ws.Rows("40:40").Insert(-4121)

# Populate the new row 40 with this week's entry.
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44964
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112028
$ws.Range("G40").Value = "Sandia"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 380
$ws.Range("L40").Value = 400
$ws.Range("M40").Value = 390
$ws.Range("N40").Value = "$/kilo (volumen en unidades)"
$ws.Range("O40").Value = "Perú"
$ws.Range("P40").Value = 390
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
